# Applies the "confirmations" summary sheet edit:
#  1. Remove the standalone "Summary" header row (row 32), shifting the
#     trailing "Total ..." rows up by one.
#  2. Re-label each category's sub-rows ("New nominations", "Confirmed",
#     "Unconfirmed", "Withdrawn", "Returned to White House") to be prefixed
#     with their parent category name (e.g. "     Civilian, New nominations").
#  3. Re-label the final "Total ..." summary rows with their new wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the "Summary" label-only row - everything below shifts up one row.
$ws.Rows(32).Delete()

# 2) Re-label category sub-rows with "<Category>, <Item>" text.
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Confirmed "
$ws.Range("A15").Value = "     Other Civilian, Unconfirmed "
$ws.Range("A16").Value = "     Other Civilian, Withdrawn "

$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Confirmed "
$ws.Range("A20").Value = "     Air Force, Unconfirmed "

$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("A23").Value = "     Army, Confirmed "
$ws.Range("A24").Value = "     Army, Unconfirmed "

$ws.Range("A26").Value = "     Navy, New nominations"
$ws.Range("A27").Value = "     Navy, Confirmed "
$ws.Range("A28").Value = "     Navy, Unconfirmed "

$ws.Range("A30").Value = "     Marine Corps, New nominations"
$ws.Range("A31").Value = "     Marine Corps, Confirmed "

# 3) Re-label the trailing "Total ..." rows (now rows 32-36 after the delete).
$ws.Range("A32").Value = "Total new nominations"
$ws.Range("A33").Value = "Total confirmed "
$ws.Range("A34").Value = "Total unconfirmed "
$ws.Range("A35").Value = "Total withdrawn "
$ws.Range("A36").Value = "Total returned to the White House "
